$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix existing row: user D5 "isMale" flag was false, should be true ---
$ws.Range("D5").Value = $true

# --- Row 6: add check is_deleted of user == false (new admin user) ---
$ws.Range("A6").Value = "fdsa@gmail.com"
[void]$ws.Hyperlinks.Add($ws.Range("A6"), "mailto:fdsa@gmail.com")

$ws.Range("B6").Value = "Kim Jong Un"
$ws.Range("B7").Value = "Kim Ji Jok"
$ws.Range("B8").Value = "Lee Chong Whey"

$ws.Range("E6").Value = "Số 5 Ngô Tất Tố"
$ws.Range("E7").Value = "Địa bàn Phường Nam Nha Trang"
$ws.Range("E8").Value = "Địa bàn Phường Bắc Từ Liêm"

$ws.Range("A7").Value = "aS@gmail.com"
[void]$ws.Hyperlinks.Add($ws.Range("A7"), "mailto:aS@gmail.com")

$ws.Range("A8").Value = "hi@gmail.com"
[void]$ws.Hyperlinks.Add($ws.Range("A8"), "mailto:hi@gmail.com")

# dob column (text, mirrors how "20/11/1990" was entered for row 5)
$ws.Range("C6").Value = "20/11/1990"
$ws.Range("C7").Value = "20/11/1990"
$ws.Range("C8").Value = "20/11/1990"

# isMale column
$ws.Range("D6").Value = $false
$ws.Range("D7").Value = $true
$ws.Range("D8").Value = $false

# phone_number column (stored as text, like the existing rows)
$ws.Range("F6").Value = "'0362718422"
$ws.Range("F7").Value = "'0123321123"
$ws.Range("F8").Value = "'0123456789"

# Row 9: blank row that still carries the hyperlink / phone-number formatting
# (mirrors the trailing formatted-but-empty row in the source workbook)
$ws.Range("A6").Copy()
$ws.Range("A9").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("F6").Copy()
$ws.Range("F9").PasteSpecial(-4122)  # xlPasteFormats

$ws.Application.CutCopyMode = $false

# Final selection, matching the author's last cursor position
[void]$ws.Range("I9").Select()
